# Weekly update: insert a new "Alcachofa" price observation as row 8,
# pushing the existing rows 8-18 down to 9-19 (dimension grows to A1:R19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8; this shifts rows 8..18 down
# to 9..19 (carrying their formatting/styles with them), matching the
# row-by-row shift seen in the diff.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly observation.
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44484
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112013
$ws.Cells.Item(8, 7).Value = "Alcachofa"
$ws.Cells.Item(8, 8).Value = "Madrigal"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 120
$ws.Cells.Item(8, 11).Value = 11000
$ws.Cells.Item(8, 12).Value = 12000
$ws.Cells.Item(8, 13).Value = 11500
$ws.Cells.Item(8, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(8, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(8, 16).Value = 288
$ws.Cells.Item(8, 17).Value = 40
$ws.Cells.Item(8, 18).Value = "Hortaliza"
